$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52/53 cleanup: two cells ("Manage permissions within the current
#     user model" and "Integrate permissions into existing") move out of
#     C52/C53 down into new rows 57/58, and the stray "Loh-dant" note in
#     D53 is deleted outright.
$ws.Range("C52").ClearContents()
$ws.Range("C53").ClearContents()
$ws.Range("D53").ClearContents()

# --- New rows 57 & 58 carrying the text that used to live in C52/C53.
$ws.Range("A57").Value = "Manage permissions within the current user model"
$ws.Range("A58").Value = "Integrate permissions into existing"

# --- Formatting: A52 gets a plain Calibri font (no fill); B52 and the
#     existing B54 ("Notice Board") get an orange font on a yellow fill.
$a52 = $ws.Range("A52")
$a52.Font.Name = "Calibri"

$b52 = $ws.Range("B52")
$b52.Interior.Color = 65535
$b52.Font.Color = 49407

$b54 = $ws.Range("B54")
$b54.Interior.Color = 65535
$b54.Font.Color = 49407

# --- Selection moves to A52 (was D58).
[void]$ws.Range("A52").Select()

Write-Host "Applied Planning.xlsx edits"
